# Updated cryptos list values (mirrors upstream data refresh commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.093.82"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.482.12"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.93"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.34"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.481.45"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.38"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.032.85"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.511.79"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.97"
$ws.Range("E19").Value = "  -5.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.94"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.53"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.608.62"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0898"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.68"
$ws.Range("E32").Value = "  -5.22%  "
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  -3.78%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.96"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  -6.74%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.22"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("E40").Value = "  -6.01%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.80"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.79"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.72"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("E50").Value = "  -6.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0729"
$ws.Range("E51").Value = "  -1.00%  "

Write-Host "Updated cryptos list values"
